$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 1.7
$ws.Range("H3").Value = 5.6
$ws.Range("I3").Value = 6.4
$ws.Range("J3").Value = 3.75
$ws.Range("Q3").Value = 1.91
